{"js": "// The evaluation section originally read:\n//   \"...para o experimento com SVM, para regress\u00e3o seria avaliado Acur\u00e1cia...\"\n// and should become:\n//   \"...para o experimento com SVM, e para todos tamb\u00e9m seria avaliado Acur\u00e1cia...\"\n// i.e. the phrase \"para regress\u00e3o\" is replaced with \"e para todos tamb\u00e9m\".\n\nconst body = context.document.body;\n\nconst results = body.search(\"para regress\u00e3o\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase 'para regress\u00e3o' not found in document body.\");\n}\n\nresults.items[0].insertText(\"e para todos tamb\u00e9m\", \"Replace\");\nawait context.sync();\n", "ps1": "# The evaluation section originally read:\n#   \"...para o experimento com SVM, para regress\u00e3o seria avaliado Acur\u00e1cia...\"\n# and should become:\n#   \"...para o experimento com SVM, e para todos tamb\u00e9m seria avaliado Acur\u00e1cia...\"\n# i.e. the phrase \"para regress\u00e3o\" is replaced with \"e para todos tamb\u00e9m\".\n\n$d = $word.ActiveDocument\n$range = $d.Content\n\n$find = $range.Find\n$find.Text = \"para regress\u00e3o\"\n$find.Replacement.Text = \"e para todos tamb\u00e9m\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdFindContinue = 1, wdReplaceOne = 1 (only one occurrence exists/is targeted)\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1) | Out-Null\n"}
